$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("57:57").Insert()
